# The underlying data records for rows 2, 4 and 5 get rotated:
#   new row2 = old row5 data
#   new row4 = old row2 data
#   new row5 = old row4 data
# (rows are re-sorted observations; columns A,B,D,E,F,G,H,Q,R,Y,AA change,
#  the rest of each row - C,I,P,S,T,U,V,W,AC,AD,AE,AG,AT,AW,AX,AY - is
#  identical across these three rows already, so no edit is needed there).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","Q","R","Y","AA")

# Capture current (pre-edit) values for the three affected rows before
# overwriting anything, since the rotation reads from one row to write
# into another.
$row2 = @{}
$row4 = @{}
$row5 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range($col + "2").Value()
    $row4[$col] = $ws.Range($col + "4").Value()
    $row5[$col] = $ws.Range($col + "5").Value()
}

# The two date-like text columns (Y, AA) must stay plain text - Excel would
# otherwise reinterpret a "yyyy-mm-dd" string as a real date serial number
# when it is assigned. Force text formatting on those cells up front.
$dateCols = @("Y","AA")
foreach ($col in $dateCols) {
    $ws.Range($col + "2").NumberFormat = "@"
    $ws.Range($col + "4").NumberFormat = "@"
    $ws.Range($col + "5").NumberFormat = "@"
}

# Row 2 <- old Row 5
foreach ($col in $cols) {
    $ws.Range($col + "2").Value = $row5[$col]
}

# Row 4 <- old Row 2
foreach ($col in $cols) {
    $ws.Range($col + "4").Value = $row2[$col]
}

# Row 5 <- old Row 4
foreach ($col in $cols) {
    $ws.Range($col + "5").Value = $row4[$col]
}

Write-Output "Row2 A/Y: $($ws.Range('A2').Value()) / $($ws.Range('Y2').Value())"
Write-Output "Row4 A/Y: $($ws.Range('A4').Value()) / $($ws.Range('Y4').Value())"
Write-Output "Row5 A/Y: $($ws.Range('A5').Value()) / $($ws.Range('Y5').Value())"
